$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.020338535308838
$ws.Range("B1").Value = 5.396963119506836
$ws.Range("C1").Value = 2.341854333877563
$ws.Range("D1").Value = 1.558398723602295
$ws.Range("E1").Value = 1.489251255989075
